# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Summary header fields ---
# VALOR MORA total
$ws.Range("E11").Value = 649600
# Cant. Trabajadores
$ws.Range("C13").Value = 1
# Cant. Periodos
$ws.Range("F13").Value = 14

# --- Detail table rows 16-29 ---
# The worker database was resorted so the arrears periods now run in
# ascending order (2209 .. 2507) instead of descending, the worker's
# updated basic salary (1200000) applies to every period row, and a new
# 14th period (2508) is appended as row 29 for the single remaining
# worker in arrears (JOHN JAIRO TEJERA MENDOZA / CC 9296780).

$rows = @(
    @{ Row = 16; Periodo = "2209"; Mora = 25600; Salario = 1200000 },
    @{ Row = 17; Periodo = "2210"; Mora = 48000; Salario = 1200000 },
    @{ Row = 18; Periodo = "2211"; Mora = 48000; Salario = 1200000 },
    @{ Row = 19; Periodo = "2410"; Mora = 48000; Salario = 1200000 },
    @{ Row = 20; Periodo = "2411"; Mora = 48000; Salario = 1200000 },
    @{ Row = 21; Periodo = "2412"; Mora = 48000; Salario = 1200000 },
    @{ Row = 22; Periodo = "2501"; Mora = 48000; Salario = 1200000 },
    @{ Row = 23; Periodo = "2502"; Mora = 48000; Salario = 1200000 },
    @{ Row = 24; Periodo = "2503"; Mora = 48000; Salario = 1200000 },
    @{ Row = 25; Periodo = "2504"; Mora = 48000; Salario = 1200000 },
    @{ Row = 26; Periodo = "2505"; Mora = 48000; Salario = 1200000 },
    @{ Row = 27; Periodo = "2506"; Mora = 48000; Salario = 1200000 },
    @{ Row = 28; Periodo = "2507"; Mora = 48000; Salario = 1200000 },
    @{ Row = 29; Periodo = "2508"; Mora = 48000; Salario = 1200000 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = "9296780"
    $ws.Cells.Item($row, 4).Value = "JOHN JAIRO TEJERA MENDOZA"
    $ws.Cells.Item($row, 5).Value = $r.Periodo
    $ws.Cells.Item($row, 6).Value = $r.Mora
    $ws.Cells.Item($row, 7).Value = $r.Salario
}

# Column D ("Nombre Trabajador") is bestFit; its autofit width shrinks now
# that the longer "WILLIAM JOSE PACHECO CASTANO" no longer appears.
$ws.Columns("D:D").ColumnWidth = 27.6
